$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "duplicate_image_filename" column (E) is populated with "NA" for the
# data rows (2 through 21) that previously left it blank.
$ws.Range("E2:E21").Value = "NA"
